$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the commit diff.
# For columns B and C (text) we set directly.
# For columns D and E we force text format ("@") before assignment so that
# numeric-looking strings (e.g. "1.003") are not auto-converted to numbers,
# preserving the original inline-string / text semantics of the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.150.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.79"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5239"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06356"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.56"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07697"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.750.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.06%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.617"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.884.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5616"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8201"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.46"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.146.13"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.656"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.53"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.35"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.957"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1196"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.269"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.97"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05454"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.465"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.370"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.566"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9547"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.780"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.402"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5678"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01586"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.879"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8321"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.028.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.29"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.794.81"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.75"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4344"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05192"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.88%  "
